$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 updates
$ws.Range("G8").Value = 3.3
$ws.Range("I8").Value = 2.35
$ws.Range("L8").Value = 3.2
$ws.Range("M8").Value = 1.11
$ws.Range("N8").Value = 6.5
$ws.Range("O8").Value = 1.5
$ws.Range("P8").Value = 2.5
$ws.Range("Q8").Value = 2.6
$ws.Range("R8").Value = 1.48
$ws.Range("S8").Value = 1.57
$ws.Range("T8").Value = 2.25
$ws.Range("AA8").Value = 34
$ws.Range("AC8").Value = 6.5
$ws.Range("AE8").Value = 17
$ws.Range("AH8").Value = 6.5
$ws.Range("AK8").Value = 23
$ws.Range("AL8").Value = 23
$ws.Range("AT8").Value = 2.25
$ws.Range("AX8").Value = 15

# Row 23 updates
$ws.Range("M23").Value = 1.05
$ws.Range("N23").Value = 11
$ws.Range("Q23").Value = 1.88
$ws.Range("R23").Value = 1.98
$ws.Range("BD23").Value = 126
